# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (per-fund holding detail) right before the
# "总计" (totals) sheet, and inserts a matching 2022-Q1 summary row at the
# top of "总计"'s data table (pushing the older quarters down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet, positioned right before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# The handle used above for positioning tracks the tab *position*, not a
# stable sheet identity - after inserting + renaming, re-resolve "总计" by
# name again so later writes land on the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# A never-written, plain "General" cell we use purely as a formats-only
# copy source, to (re)stamp cells back to the workbook's default style.
$blank = $q1.Range("Z100")

# Header row (same look as the other quarterly sheets: bold/boxed style)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$totalSheet.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Fund-holding detail rows
$rows = @(
    @{ idx = 0; code = "110002"; name = "易方达策略成长混合";                           size = "12.15"; pos = "88.76"; ratio = "3.05"; value = "0.3706"; rank = 8 },
    @{ idx = 1; code = "112002"; name = "易方达策略成长二号混合";                       size = "10.72"; pos = "87.99"; ratio = "3.02"; value = "0.3237"; rank = 9 },
    @{ idx = 2; code = "159851"; name = "华宝中证金融科技主题ETF";                       size = "3.16";  pos = "98.58"; ratio = "3.02"; value = "0.0954"; rank = 8 },
    @{ idx = 3; code = "516100"; name = "华夏中证金融科技主题交易型开放式指数证券投资基金"; size = "0.68";  pos = "96.91"; ratio = "3.03"; value = "0.0206"; rank = 8 },
    @{ idx = 4; code = "004250"; name = "银河量化优选混合";                             size = "0.39";  pos = "80.03"; ratio = "2.04"; value = "0.0080"; rank = 3 }
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row.idx

    # Text-valued columns: force text storage (so "12.15" etc. stay strings,
    # not numbers), then strip the resulting NumberFormat-driven style back
    # to the sheet's default so the cells end up unstyled, same as the rest
    # of the data rows.
    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row.code
    $q1.Cells.Item($r, 3).NumberFormat = "@"
    $q1.Cells.Item($r, 3).Value = $row.name
    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row.size
    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row.pos
    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row.ratio
    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row.value
    $q1.Cells.Item($r, 8).Value = $row.rank

    $blank.Copy()
    $q1.Range($q1.Cells.Item($r, 2), $q1.Cells.Item($r, 7)).PasteSpecial(-4122)

    $r = $r + 1
}

# Column A (row-index column) uses the same bold/boxed style as the header
$totalSheet.Range("B1").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

$blank.Clear()

# ---------------------------------------------------------------------
# 2. Insert a new top data row in "总计" for 2022-Q1, push the rest down
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 5
$totalSheet.Cells.Item(2, 4).Value = 0.82

# Row-insert can drag in stray formatting - restamp to match the sheet's
# existing data rows: column A keeps the bold/boxed style, B:D stay plain.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("B3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

# Renumber the index column (A) for the rows that shifted down
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3

$excel.CutCopyMode = $false
